$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first 7 data rows (rows 2-8), which shifts everything up.
$ws.Rows("2:8").Delete()

# New data to append after the shift (now the last existing data row is 15).
$newData = @(
    @(-3.313027620315552, 2.173916578292847, -4.959309577941895),
    @(0.2906191349029541, -1.299157619476318, 1.514182209968567),
    @(1.71515691280365, 0.9166033267974854, 0.9940304756164552),
    @(-1.992336988449097, -0.5707008838653564, 0.5294674634933472),
    @(-0.0847575515508651, -0.578489363193512, 2.060448408126831),
    @(1.970193147659302, -0.4350887835025787, 1.578170418739319)
)

$startRow = 16
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newData[$i][0]
    $ws.Cells.Item($row, 2).Value = $newData[$i][1]
    $ws.Cells.Item($row, 3).Value = $newData[$i][2]
}

$wb.Save()
